$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "252.23"
$cell.Style = "Normal"

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "22.69"
$cell.Style = "Normal"

# Row 4
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "5.419"
$cell.Style = "Normal"

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "0.05714"
$cell.Style = "Normal"

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "3.417"
$cell.Style = "Normal"

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "6.365"
$cell.Style = "Normal"

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.9400"
$cell.Style = "Normal"

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.1437"
$cell.Style = "Normal"

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07481"
$cell.Style = "Normal"

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.03145"
$cell.Style = "Normal"

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.03077"
$cell.Style = "Normal"

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "3.717"
$cell.Style = "Normal"

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.001584"
$cell.Style = "Normal"

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.04776"
$cell.Style = "Normal"

# Row 18
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.0005786"
$cell.Style = "Normal"
$ws.Range("E18").Value = "17OneONE"

# Row 19
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.006422"
$cell.Style = "Normal"
$ws.Range("E19").Value = "18TigerCashTCH"

# Row 20
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.005044"
$cell.Style = "Normal"
$ws.Range("E20").Value = "19HotbitTokenHTB"

# Row 21
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.001027"
$cell.Style = "Normal"
$ws.Range("E21").Value = "20BitKanKAN"

# Row 22
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.0001499"
$cell.Style = "Normal"
$ws.Range("E22").Value = "21NitroExNTX"

# Row 23
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "3.710"
$cell.Style = "Normal"
$ws.Range("E23").Value = "22LEOLEO"

# Row 24
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.179"
$cell.Style = "Normal"
$ws.Range("E24").Value = "23BTSETokenBTSE"

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.1307"
$cell.Style = "Normal"

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.0002998"
$cell.Style = "Normal"

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.04027"
$cell.Style = "Normal"

# Row 41
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.1072"
$cell.Style = "Normal"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# Row 42
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.002709"
$cell.Style = "Normal"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.002937"
$cell.Style = "Normal"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.008066"
$cell.Style = "Normal"

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.00005756"
$cell.Style = "Normal"

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.00000000749"
$cell.Style = "Normal"

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.4996"
$cell.Style = "Normal"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.00002098"
$cell.Style = "Normal"

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.01009"
$cell.Style = "Normal"
